$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.177251577377319
$ws.Range("B1").Value = 2.418758630752563
$ws.Range("D1").Value = 2.334059715270996
$ws.Range("E1").Value = 1.2027268409729
